# Fruta / hortaliza, semanal
# The weekly data refresh re-ordered the 5 price observation rows (rows 2-6)
# of the sheet. Concretely, new row N gets the values that used to live in a
# different source row, per the mapping below (new row -> old/source row):
#   2 <- 3
#   3 <- 4
#   4 <- 6
#   5 <- 2
#   6 <- 5
# Capture the "before" values first, then write them back in the new order,
# so this works regardless of execution order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually vary between rows in this block (per the diff):
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion),
# P (Precio $/Kg), Q (Kg o Unidades)
$cols = @("D", "J", "K", "L", "M", "N", "P", "Q")

# Snapshot current ("before") values for rows 2..6.
# NOTE: use Value2 (not Value) to read - Value is a parameterized COM
# property and reading it directly yields the property descriptor rather
# than the underlying scalar in this environment.
$before = @{}
foreach ($r in 2..6) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowVals
}

# new row -> source row mapping
$map = @{ 2 = 3; 3 = 4; 4 = 6; 5 = 2; 6 = 5 }

foreach ($newRow in 2..6) {
    $srcRow = $map[$newRow]
    $srcVals = $before[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $srcVals[$c]
    }
}
